$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28: A/Q/R rotate in from old row 32
$ws.Range("A28").Value = 111363037
$ws.Range("Q28").Value = 590869.45745891
$ws.Range("R28").Value = 6985763.548076616

# Row 30: A/Q/R rotate in from old row 28
$ws.Range("A30").Value = 111363036
$ws.Range("Q30").Value = 590784.7419292277
$ws.Range("R30").Value = 6985796.657093059

# Row 31 and Row 33 fully swap their species-related data (B,D,E,F,G,H) plus own A/Q/R
$ws.Range("A31").Value = 111363040
$ws.Range("B31").Value = 89405
$ws.Range("D31").Value = "NT"
$ws.Range("E31").Value = 1202
$ws.Range("F31").Value = "Ullticka"
$ws.Range("G31").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H31").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q31").Value = 590846.7337154552
$ws.Range("R31").Value = 6985549.62329052

# Row 32: A/Q/R rotate in from old row 31
$ws.Range("A32").Value = 111363038
$ws.Range("Q32").Value = 590741.4330251076
$ws.Range("R32").Value = 6985602.172479901

# Row 33: gets the species data that used to be in row 31
$ws.Range("A33").Value = 111363042
$ws.Range("B33").Value = 96348
$ws.Range("D33").Value = "VU"
$ws.Range("E33").Value = 220787
$ws.Range("F33").Value = "Knärot"
$ws.Range("G33").Value = "Goodyera repens"
$ws.Range("H33").Value = "(L.) R. Br."
$ws.Range("Q33").Value = 590767.8557530388
$ws.Range("R33").Value = 6985781.634191129
